# Fix lỗi tạo dữ liệu excel
# - Header E1 changes from "Trạng thái" to "Check QC"
# - Data cell C2 keeps its text "PSSM Tẩy Hồng 009 CYCLAMEN (XK/CR)"
# - Active selection moves from F9 to F17

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Check QC"
$ws.Range("C2").Value = "PSSM Tẩy Hồng 009 CYCLAMEN (XK/CR)"

$ws.Range("F17").Select()
